$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (25) down to the new row (26)
$ws.Range("A25:B25").Copy()
$ws.Range("A26:B26").PasteSpecial(-4122) # xlPasteFormats

# Fill in the new timelog entry
$ws.Range("A26").Value = "3/10, 4 hrs"
$ws.Range("B26").Value = "Working on adding govt measures to graphs, documenting, fixing errors, prepping final product"

# Match the row height used by similar longer entries
$ws.Rows.Item(26).RowHeight = 55.2

# Move the active selection to where the user would type the next entry
$null = $ws.Range("B27").Select()
